# Insert a new weekly price record before the current row 38, shifting the
# existing rows 38:66 down to 39:67 (dimension grows from A1:R66 to A1:R67).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(38).Insert()

# Populate the newly-inserted row 38 with the new "Espinaca" price record.
$ws.Cells.Item(38, 1).Value = 4
$ws.Cells.Item(38, 2).Value = 'Feria Lagunitas de Puerto Montt'
$ws.Cells.Item(38, 3).Value = 'Los Lagos'
$ws.Cells.Item(38, 4).Value = 45086
$ws.Cells.Item(38, 5).Value = 10
$ws.Cells.Item(38, 6).Value = 100112012
$ws.Cells.Item(38, 7).Value = 'Espinaca'
$ws.Cells.Item(38, 8).Value = 'Sin especificar'
$ws.Cells.Item(38, 9).Value = 'Primera'
$ws.Cells.Item(38, 10).Value = 25
$ws.Cells.Item(38, 11).Value = 12000
$ws.Cells.Item(38, 12).Value = 12000
$ws.Cells.Item(38, 13).Value = 12000
$ws.Cells.Item(38, 14).Value = '$/cuna 10 kilos'
$ws.Cells.Item(38, 15).Value = 'Región Metropolitana'
$ws.Cells.Item(38, 16).Value = 1200
$ws.Cells.Item(38, 17).Value = 10
$ws.Cells.Item(38, 18).Value = 'Hortaliza'
